# Add data for 2021-09-30
# The "Through" date advances from 2021-09-21 to 2021-09-22, and the
# underlying carjacking counts for a handful of neighborhood/month cells
# are incremented (or newly populated) to reflect the newly added record(s).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet title and header label to reflect the new "through" date.
$ws.Name = "Through 2021-09-22"
$ws.Range("B1").Value = "September 2021 (through September 22)"

# Updated / newly populated monthly counts per neighborhood.
$ws.Range("B3").Value = 8
$ws.Range("B5").Value = 9
$ws.Range("K5").Value = 10
$ws.Range("AL5").Value = 5
$ws.Range("K6").Value = 4
$ws.Range("B11").Value = 4
$ws.Range("T11").Value = 1
$ws.Range("B13").Value = 6
$ws.Range("T18").Value = 1
$ws.Range("T19").Value = 3
$ws.Range("T22").Value = 3
$ws.Range("B23").Value = 4
$ws.Range("T39").Value = 2
$ws.Range("AL53").Value = 2
$ws.Range("AC55").Value = 4
$ws.Range("B60").Value = 2
$ws.Range("AC66").Value = 1
